# Further experiments with HFD: k variation, half-sampling rate
#
# Inserts a new results row (new row 18: "HFD + NuSVM (linear kernel)",
# nu = 0.25007, HFD over 375-500) into the results table on Sheet1,
# pushing the previous rows 18-22 down to 19-23, and updates the
# worksheet selection to match where the author ended up (F18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 18:22 down one row, opening up row 18 for the new entry.
$ws.Rows.Item(18).Insert()

# New row 18 data.
$ws.Cells.Item(18, 1).Value = "HFD + NuSVM (linear kernel)"
$ws.Cells.Item(18, 2).Value = 0.8349
$ws.Cells.Item(18, 3).Value = "17/19"
$ws.Cells.Item(18, 4).Value = "RH"
$ws.Cells.Item(18, 5).Value = 6
$ws.Cells.Item(18, 6).Value = "nu = 0.25007, HFD over 375-500"

# Match the author's final on-screen selection.
$ws.Range("F18").Select()
